# Update "想去人数" (number of people interested) counts for a few
# exhibition entries, in both the "展览" sheet and the combined
# "全部类型" sheet.

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions) sheet
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 513   # 南宁·布谷鸟动漫展5th: 512 -> 513
$wsExpo.Range("F3").Value = 6117  # 南宁·2024良牙动漫秋季盛典（秋典）: 6110 -> 6117
$wsExpo.Range("F6").Value = 119   # 南宁·花海演绎二次元水上派对: 118 -> 119
$wsExpo.Range("F9").Value = 559   # 南宁·万圣漫控嘉年华10: 558 -> 559

# 全部类型 (All types, combined) sheet
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 513    # 南宁·布谷鸟动漫展5th: 512 -> 513
$wsAll.Range("F3").Value = 6117   # 南宁·2024良牙动漫秋季盛典（秋典）: 6110 -> 6117
$wsAll.Range("F7").Value = 119    # 南宁·花海演绎二次元水上派对: 118 -> 119
$wsAll.Range("F11").Value = 559   # 南宁·万圣漫控嘉年华10: 558 -> 559
